$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shelf entries (weight and content filters sample data)
$ws.Range("A14").Value = "E1"
$ws.Range("B14").Value = 40
$ws.Range("C14").Value = "İÇECEK"
$ws.Range("D14").Value = "AĞIR"

$ws.Range("A15").Value = "E2"
$ws.Range("B15").Value = 15
$ws.Range("C15").Value = "YİYECEK"
$ws.Range("D15").Value = "HAFİF"

$ws.Range("A16").Value = "E3"
$ws.Range("B16").Value = 60
$ws.Range("C16").Value = "İÇECEK"
$ws.Range("D16").Value = "AĞIR"

$ws.Range("D17").Select() | Out-Null
